$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($cellRef, $value)
    $c = $ws.Range($cellRef)
    $c.Value = "'" + $value
    $c.Style = "Normal"
}

Set-TextValue "D2" "27.149.26"
Set-TextValue "E2" "  -2.95%  "
Set-TextValue "D3" "1.711.59"
Set-TextValue "E3" "  -3.47%  "
Set-TextValue "D4" "1.002"
Set-TextValue "E4" "  +0.10%  "
Set-TextValue "D5" "308.56"
Set-TextValue "E5" "  -6.03%  "
Set-TextValue "E6" "  +0.14%  "
Set-TextValue "D7" "0.4756"
Set-TextValue "E7" "  +5.75%  "
Set-TextValue "D8" "0.3438"
Set-TextValue "E8" "  -3.41%  "
Set-TextValue "D9" "42.12"
Set-TextValue "E9" "  +0.43%  "
Set-TextValue "D10" "0.07275"
Set-TextValue "E10" "  -2.22%  "
Set-TextValue "D11" "1.042"
Set-TextValue "E11" "  -5.93%  "
Set-TextValue "D12" "1.002"
Set-TextValue "E12" "  +0.15%  "
Set-TextValue "D13" "19.84"
Set-TextValue "E13" "  -5.58%  "
Set-TextValue "D14" "5.858"
Set-TextValue "E14" "  -3.21%  "
Set-TextValue "D15" "1.710.75"
Set-TextValue "E15" "  -3.48%  "
Set-TextValue "D16" "6.844"
Set-TextValue "E16" "  -5.74%  "
Set-TextValue "D17" "88.74"
Set-TextValue "E17" "  -5.48%  "
Set-TextValue "E18" "  -2.25%  "
Set-TextValue "D19" "0.06356"
Set-TextValue "E19" "  -1.38%  "
Set-TextValue "D20" "1.001"
Set-TextValue "E20" "  +0.18%  "
Set-TextValue "D21" "16.49"
Set-TextValue "E21" "  -4.03%  "
Set-TextValue "D22" "5.606"
Set-TextValue "E22" "  -3.21%  "
Set-TextValue "D23" "27.183.32"
Set-TextValue "E23" "  -2.91%  "
Set-TextValue "D24" "10.80"
Set-TextValue "E24" "  -4.62%  "
Set-TextValue "D25" "2.091"
Set-TextValue "E25" "  -1.76%  "
Set-TextValue "D26" "153.48"
Set-TextValue "E26" "  -5.18%  "
Set-TextValue "D27" "19.67"
Set-TextValue "E27" "  -3.61%  "
Set-TextValue "D28" "1.906.05"
Set-TextValue "E28" "  -3.53%  "
Set-TextValue "D29" "2.080"
Set-TextValue "E29" "  -4.08%  "
Set-TextValue "D30" "119.97"
Set-TextValue "E30" "  -3.87%  "
Set-TextValue "D31" "1.015"
Set-TextValue "E31" "  -8.70%  "
Set-TextValue "D32" "0.09252"
Set-TextValue "E32" "  +0.30%  "
Set-TextValue "D33" "3.591"
Set-TextValue "E33" "  -2.73%  "
Set-TextValue "D34" "5.294"
Set-TextValue "E34" "  -7.32%  "
Set-TextValue "D35" "0.02192"
Set-TextValue "E35" "  -4.36%  "
Set-TextValue "D36" "0.05893"
Set-TextValue "E36" "  -5.10%  "
Set-TextValue "D37" "11.04"
Set-TextValue "E37" "  -7.05%  "
Set-TextValue "D38" "0.2010"
Set-TextValue "E38" "  -4.86%  "
Set-TextValue "E39" "  -4.78%  "
Set-TextValue "D40" "1.412"
Set-TextValue "E40" "  +1.02%  "
Set-TextValue "D41" "1.001"
Set-TextValue "E41" "  +0.22%  "
Set-TextValue "D42" "0.5920"
Set-TextValue "E42" "  -6.35%  "
Set-TextValue "E43" "  -6.13%  "
Set-TextValue "D44" "7.479"
Set-TextValue "E44" "  -5.49%  "
Set-TextValue "D45" "12.65"
Set-TextValue "E45" "  -4.76%  "
Set-TextValue "D46" "3.566"
Set-TextValue "E46" "  -4.90%  "
Set-TextValue "D47" "0.5620"
Set-TextValue "E47" "  -4.67%  "
Set-TextValue "D48" "118.24"
Set-TextValue "E48" "  -3.71%  "
Set-TextValue "D49" "1.839"
Set-TextValue "E49" "  -6.31%  "
Set-TextValue "D50" "0.06634"
Set-TextValue "E50" "  -3.77%  "
Set-TextValue "D51" "1.085"
Set-TextValue "E51" "  -5.10%  "
